$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per row, matching the diff
$updates = @{
    2  = -3
    3  = -7
    4  = -5
    5  = -6
    6  = -2
    7  = -1
    8  = 2
    9  = -2
    10 = -2
    12 = -2
    14 = -5
    15 = 4
    16 = 0
    17 = -3
    19 = 1
    21 = -6
    22 = 3
    23 = -4
    24 = -5
    25 = 3
    26 = -3
    27 = -2
    28 = -5
    29 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
